# Update column F (dSF) values for specific rows, per repull/mean recalculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2
$ws.Range("F3").Value = -4
$ws.Range("F5").Value = -6
$ws.Range("F6").Value = -3
$ws.Range("F12").Value = 1
$ws.Range("F13").Value = -11
$ws.Range("F15").Value = -7
